# Update countries & provincias Spain
# - Refresh the "last updated" timestamp string
# - Refresh country case/recovered/death figures for several rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Updated timestamp (was "...a las 13:05")
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 13:35"

# India (row 13): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B13").Value = 153237
$ws.Range("C13").Value = 2444
$ws.Range("D13").Value = 64733
$ws.Range("E13").Value = 84139

# Austria (row 44): Casos activos, Recuperados, Muertes hoy, Muertes
$ws.Range("D44").Value = 15228
$ws.Range("E44").Value = 684
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 645

# Kazajistan (row 54): Casos activos, Recuperados
$ws.Range("D54").Value = 4746
$ws.Range("E54").Value = 4521

# Republica de Macedonia (row 90): Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes
$ws.Range("B90").Value = 2039
$ws.Range("C90").Value = 25
$ws.Range("D90").Value = 1470
$ws.Range("E90").Value = 450
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 119

# San Marino (row 132): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B132").Value = 667
$ws.Range("C132").Value = 1
$ws.Range("D132").Value = 294
$ws.Range("E132").Value = 331
